$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point difference on A86 (precision correction)
$ws.Cells.Item(86, 1).Value = 44399.76833309606

# Add new row 87 of data
$ws.Cells.Item(87, 1).Value = 44400.76775271606
$ws.Cells.Item(87, 2).Value = 80385
$ws.Cells.Item(87, 3).Value = 67759
$ws.Cells.Item(87, 4).Value = 3737
$ws.Cells.Item(87, 5).Value = 2217
$ws.Cells.Item(87, 6).Value = 1597
$ws.Cells.Item(87, 7).Value = 21066
$ws.Cells.Item(87, 8).Value = 1664
$ws.Cells.Item(87, 9).Value = 894
$ws.Cells.Item(87, 10).Value = 207

# Ensure A87 has same number format as the rest of column A (date format)
$ws.Cells.Item(87, 1).NumberFormat = $ws.Cells.Item(86, 1).NumberFormat
